# Apply the LOT2022.xlsx content realignment:
# - Several label/value cells in rows 10-24 were reshuffled up (a handful of
#   large descriptive paragraphs were removed from the shared-string table,
#   and the surviving labels/values shifted to fill the gaps).
# - Row 25 (the last "LOQ4086 ..." requirement row) no longer exists; the
#   workbook now ends at row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("B10").Value = '6007846 - Júlio César dos Santos'
$ws.Range("C10").Value = '6007846 - Júlio César dos Santos'

# Row 13
$ws.Rows.Item(13).RowHeight = 60
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'

# Row 14
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Introduction to modeling and simulation of bioprocesses; study of problems of the industry of bioprocesses related to the construction and solution of phenomenological models: computational software and algebraic equations; mathematical modeling and simulation of fermentative processes; constructing and solving models: differential equations; adjustment of parameters and bioprocesses optimization; use of process simulators applied to biotechnology.'
$ws.Range("C14").Value = 'Introduction to modeling and simulation of bioprocesses; study of problems of the industry of bioprocesses related to the construction and solution of phenomenological models: computational software and algebraic equations; mathematical modeling and simulation of fermentative processes; constructing and solving models: differential equations; adjustment of parameters and bioprocesses optimization; use of process simulators applied to biotechnology.'

# Row 15
$ws.Rows.Item(15).RowHeight = 120
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2018'
$ws.Range("C15").Value = '01/01/2018'

# Row 16
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1. Introduction to modeling and simulation of bioprocesses.1.1. Definition of mathematical model.1.2. Concepts of dependent and independent variables of a system.1.3. Definition and classification of control volumes.2. Study of problems of the industry of bioprocesses related to the construction and solution of phenomenological models: computational software and algebraic equations.2.1. Introduction to computational software/packages used to solving mathematical models.2.2. Solving of problems using systems of linear equations.2.3. Solving of problems using non-linear equations.2.4. Solving of problems using systems of non-linear equations.3. Mathematical modeling and simulation of fermentative processes3.1. Objectives3.2. Differences between chemical and fermentative processes3.3. Interactions between the microbial population and the culture medium.3.4. Construction and classification of mathematical models for fermentative processes.3.5. Kinetic models of cellular growth, substrate consumption and formation of products in fermentative processes.3.6. Modeling of fermentative process in reactors: batch, continuous, continuous with cells recycle, fed-batch and tubular.4. Constructing and solving models: differential equations.5. Adjustment of parameters and bioprocesses optimization.6. Use of process simulators applied to biotechnology.6.1. Process design aided by simulation software6.2. Classification of process simulation software 6.3. Synthesis and analysis of process6.4. Process flowsheeting: concepts and limitation, convergence6.5. Application examples.'
$ws.Range("C16").Value = '1. Introduction to modeling and simulation of bioprocesses.1.1. Definition of mathematical model.1.2. Concepts of dependent and independent variables of a system.1.3. Definition and classification of control volumes.2. Study of problems of the industry of bioprocesses related to the construction and solution of phenomenological models: computational software and algebraic equations.2.1. Introduction to computational software/packages used to solving mathematical models.2.2. Solving of problems using systems of linear equations.2.3. Solving of problems using non-linear equations.2.4. Solving of problems using systems of non-linear equations.3. Mathematical modeling and simulation of fermentative processes3.1. Objectives3.2. Differences between chemical and fermentative processes3.3. Interactions between the microbial population and the culture medium.3.4. Construction and classification of mathematical models for fermentative processes.3.5. Kinetic models of cellular growth, substrate consumption and formation of products in fermentative processes.3.6. Modeling of fermentative process in reactors: batch, continuous, continuous with cells recycle, fed-batch and tubular.4. Constructing and solving models: differential equations.5. Adjustment of parameters and bioprocesses optimization.6. Use of process simulators applied to biotechnology.6.1. Process design aided by simulation software6.2. Classification of process simulation software 6.3. Synthesis and analysis of process6.4. Process flowsheeting: concepts and limitation, convergence6.5. Application examples.'

# Row 17
$ws.Rows.Item(17).AutoFit()
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

# Row 18
$ws.Rows.Item(18).RowHeight = 60
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '6007846 - Júlio César dos Santos'
$ws.Range("C18").Value = '6007846 - Júlio César dos Santos'

# Row 19
$ws.Range("A19").Value = 'Critério:'

# Row 20
$ws.Range("A20").Value = 'Norma de recuperação:'

# Row 21
$ws.Rows.Item(21).RowHeight = 120
$ws.Range("A21").Value = 'Bibliografia:'

# Row 22
$ws.Rows.Item(22).AutoFit()
$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

# Row 23
$ws.Rows.Item(23).RowHeight = 30
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = 'LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)
'
$ws.Range("C23").Value = 'LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)
'

# Row 24
$ws.Range("B24").Value = 'LOQ4086 -  Operações Unitárias II  (Requisito fraco)
'
$ws.Range("C24").Value = 'LOQ4086 -  Operações Unitárias II  (Requisito fraco)
'

# The former row 25 ("LOQ4086 - Operacoes Unitarias II ...") is gone; its
# content already landed on row 24 above, so just drop the now-empty last row.
$ws.Rows.Item(25).Delete()
